$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 417, shifting existing rows 417:453 down to 418:454
$ws.Rows.Item(417).Insert()

# Populate the newly inserted row 417 with the new weekly price record
$ws.Range("A417").Value = 5
$ws.Range("B417").Value = "Macroferia Regional de Talca"
$ws.Range("C417").Value = "Maule"
$ws.Range("D417").Value = 45106
$ws.Range("E417").Value = 7
$ws.Range("F417").Value = 100112009
$ws.Range("G417").Value = "Acelga"
$ws.Range("H417").Value = "Sin especificar"
$ws.Range("I417").Value = "Primera"
$ws.Range("J417").Value = 500
$ws.Range("K417").Value = 2000
$ws.Range("L417").Value = 2000
$ws.Range("M417").Value = 2000
$ws.Range("N417").Value = "$/docena de atados (4 kilos)"
$ws.Range("O417").Value = "Región del Maule"
$ws.Range("P417").Value = 500
$ws.Range("Q417").Value = 4
$ws.Range("R417").Value = "Hortaliza"
